$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("F8").Value = "Desenho Técnico"
$ws.Range("E9").Value = "Desenho Técnico"
$ws.Range("D12").Value = "-"
